# Update the sample data on the "Ad_Feedback" sheet: replace the four
# feedback rows with a new set of Ad/Question/Text values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ad_Feedback")

# Set the new "Text" values first, in the same order the new strings were
# introduced, so the underlying shared-string table is built up the same way.
$ws.Range("B3").Value = "vaping can be dangerous to immune system"
$ws.Range("B2").Value = "dishonesty in industry"
$ws.Range("B5").Value = "the harms of vaping to your lungs"
$ws.Range("B4").Value = "It's not safe to vape"

# Row 2: Text="dishonesty in industry", Ad="DD", Question="Main Message"
$ws.Range("C2").Value = "DD"
$ws.Range("D2").Value = "Main Message"

# Row 3: Text="vaping can be dangerous to immune system", Ad="DF", Question="Main Message"
$ws.Range("C3").Value = "DF"
$ws.Range("D3").Value = "Main Message"

# Row 4: Text="It's not safe to vape", Ad="ST", Question="Main Message"
$ws.Range("C4").Value = "ST"
$ws.Range("D4").Value = "Main Message"

# Row 5: Text="the harms of vaping to your lungs", Ad="DF", Question="Main Message"
$ws.Range("C5").Value = "DF"
$ws.Range("D5").Value = "Main Message"

# Make Ad_Feedback the active sheet/tab, with D7 selected (matches the
# workbook's saved view state).
$ws.Select() | Out-Null
$ws.Range("D7").Select() | Out-Null
